# Applies the "Updated solver file and example datasets" commit:
#  - Variable_data: new row 3 (RadS example) + its formatting
#  - Variable_data: selection moves to I3
#  - Variable_data: D2:D48 list validation becomes a "normal" data validation
#    (in addition to the existing x14 extension entry)
#  - Distributions: sheetView selection moves to B38

$wb = $excel.ActiveWorkbook

# ---- Variable_data sheet ------------------------------------------------
$varSheet = $wb.Worksheets.Item("Variable_data")

# New data row (mirrors the layout/format of row 2): RadS / uniform / 0.5 / 2500000 / Example
$varSheet.Range("E2").Copy()
$varSheet.Range("E3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$varSheet.Range("A3").Value = "RadS"
$varSheet.Range("E3").Value = "uniform"
$varSheet.Range("F3").Value = 0.5
$varSheet.Range("G3").Value = 2500000
$varSheet.Range("I3").Value = "Example"

# Add the plain (non-x14) list validation for D2:D48 that points at the
# broken #REF! range (mirrors the worksheet's existing, already-broken rule).
$varSheet.Range("D2:D48").Validation.Add(3, 1, 1, "#REF!")

# Move the visible selection to I3, matching the saved sheetView state.
$varSheet.Activate()
$varSheet.Range("I3").Select()

# ---- Distributions sheet -------------------------------------------------
$distSheet = $wb.Worksheets.Item("Distributions")
$distSheet.Activate()
$distSheet.Range("B38").Select()

# Restore the originally active sheet/selection (Variable_data).
$varSheet.Activate()
